# Append: 2026-02-08 12:43 JST
#
# A new scraped work item ("初回 不動産業向け:レインズ等から...") is inserted
# into the "ランサーズ" sheet right after the header rows that precede it in
# priority-score order, pushing the previously-existing item rows down by
# one. Every "取得日時" (fetched-at) timestamp in column A is refreshed to
# the new run time, and a hyperlink is wired up for the new row's URL cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-08 12:43:54"

# Hyperlink objects don't automatically renumber/relocate when rows are
# inserted, so drop them all now and recreate them (in column order) once
# the new row is in place and every cell holds its final value.
$ws.Hyperlinks.Delete()

# Insert a new row above row 3 (the second data row); this shifts the
# previously-existing rows 3-5 down to rows 4-6, carrying their
# formatting (incl. the hyperlink style on column F) along.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new work item.
$ws.Cells.Item(3, 1).Value = $newTimestamp
$ws.Cells.Item(3, 2).Value = "初回 不動産業向け:レインズ等からの物件情報収集および社内ソフトへの自動入力ツール開発"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5487945"
$ws.Cells.Item(3, 7).Value = 135
$ws.Cells.Item(3, 8).Value = "◆ツール,開発"

# Refresh the fetched-at timestamp for every data row (2 through 6) to
# reflect this run.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Re-wire the URL hyperlink for every data row now that they're all
# settled into their final positions. The upstream export keeps its
# hyperlink relationships in their original (pre-shift) order and simply
# appends the newest link last, so rId1-4 are recreated against the same
# four legacy targets they always pointed at and the brand-new URL is
# added last as rId5, landing on whichever row is now last (F6).
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.lancers.jp/work/detail/5487791") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.lancers.jp/work/detail/5487838") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.lancers.jp/work/detail/5487908") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.lancers.jp/work/detail/5487828") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.lancers.jp/work/detail/5487945") | Out-Null

# Adding a hyperlink can fork a near-duplicate "Hyperlink" style; snap
# every URL cell back onto the sheet's single shared Hyperlink style.
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
